{"js": "// Helper: wrap a single <w:p>...</w:p> body fragment in the minimal OOXML\n// \"flat\" package that Word.Range.insertOoxml() expects.\nfunction wrapParagraphPackage(paragraphXml) {\n  return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' +\n    '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' +\n    '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' +\n    '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' +\n    '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' +\n    '</Relationships></pkg:xmlData></pkg:part>' +\n    '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' +\n    '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' +\n    paragraphXml +\n    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>';\n}\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Locate the two anchor paragraphs we need by their (unique) text content\n// instead of hard-coded indices, so the script is resilient to minor shifts.\nlet companyInfoPara = null;\nlet shortpositiesPara = null;\nlet shortPctPara = null;\nlet voorstelKoopVerkoopPara = null;\nlet overigePara = null;\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text === \"[COMPANY_INFO]\") {\n    companyInfoPara = paragraphs.items[i];\n  } else if (text === \"Shortposities\") {\n    shortpositiesPara = paragraphs.items[i];\n  } else if (text === \"[Short % of Shares Outstanding:]\") {\n    shortPctPara = paragraphs.items[i];\n  } else if (text === \"Voorstel Koop & Verkoop\") {\n    voorstelKoopVerkoopPara = paragraphs.items[i];\n  } else if (text === \"Overige\") {\n    overigePara = paragraphs.items[i];\n  }\n}\n\n// 1) Add a new \"Logo:\" paragraph (style \"Plattetekst\") right before the\n//    \"[COMPANY_INFO]\" placeholder paragraph in the \"Bedrijfsinformatie\" section.\nif (companyInfoPara) {\n  const logoPara = companyInfoPara.insertParagraph(\"Logo:\", Word.InsertLocation.before);\n  logoPara.style = \"Plattetekst\";\n  await context.sync();\n}\n\n// 2) The <w:lastRenderedPageBreak/> marker moves from the first run of the\n//    \"[Short % of Shares Outstanding:]\" paragraph up onto the \"Shortposities\"\n//    heading run (it now renders at the top of that section instead).\nif (shortpositiesPara) {\n  shortpositiesPara.getRange(\"Whole\").insertOoxml(\n    wrapParagraphPackage(\n      '<w:p><w:pPr><w:pStyle w:val=\"Kop4\"/><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr></w:pPr>' +\n      '<w:bookmarkStart w:id=\"7\" w:name=\"shortposities\"/><w:bookmarkEnd w:id=\"6\"/>' +\n      '<w:r><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr><w:lastRenderedPageBreak/><w:t>Shortposities</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\nif (shortPctPara) {\n  shortPctPara.getRange(\"Whole\").insertOoxml(\n    wrapParagraphPackage(\n      '<w:p><w:pPr><w:pStyle w:val=\"Plattetekst\"/></w:pPr>' +\n      '<w:r><w:t>[</w:t></w:r><w:r><w:t>Short % of Shares Outstanding:</w:t></w:r><w:r><w:t>]</w:t></w:r></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 3) Add a <w:lastRenderedPageBreak/> to the \"Voorstel Koop & Verkoop\" heading.\nif (voorstelKoopVerkoopPara) {\n  voorstelKoopVerkoopPara.getRange(\"Whole\").insertOoxml(\n    wrapParagraphPackage(\n      '<w:p><w:pPr><w:pStyle w:val=\"Kop4\"/></w:pPr>' +\n      '<w:bookmarkStart w:id=\"14\" w:name=\"voorstel-koop-verkoop\"/><w:bookmarkEnd w:id=\"13\"/>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:lastRenderedPageBreak/><w:t>Voorstel</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' +\n      '<w:r><w:t xml:space=\"preserve\"> Koop &amp; </w:t></w:r>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Verkoop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n\n// 4) Remove the <w:lastRenderedPageBreak/> from the \"Overige\" heading (the\n//    page break is now rendered earlier, at \"Voorstel Koop & Verkoop\").\nif (overigePara) {\n  overigePara.getRange(\"Whole\").insertOoxml(\n    wrapParagraphPackage(\n      '<w:p><w:pPr><w:pStyle w:val=\"Kop4\"/></w:pPr>' +\n      '<w:bookmarkStart w:id=\"15\" w:name=\"overige\"/><w:bookmarkEnd w:id=\"14\"/>' +\n      '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Overige</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n    ),\n    Word.InsertLocation.replace\n  );\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Wrap-ParagraphPackage {\n    param([string]$ParagraphXml)\n    return '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + `\n        '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + `\n        '<pkg:part pkg:name=\"/_rels/.rels\" pkg:contentType=\"application/vnd.openxmlformats-package.relationships+xml\">' + `\n        '<pkg:xmlData><Relationships xmlns=\"http://schemas.openxmlformats.org/package/2006/relationships\">' + `\n        '<Relationship Id=\"rId1\" Type=\"http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument\" Target=\"word/document.xml\"/>' + `\n        '</Relationships></pkg:xmlData></pkg:part>' + `\n        '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + `\n        '<pkg:xmlData><w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\"><w:body>' + `\n        $ParagraphXml + `\n        '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'\n}\n\n# Locate the anchor paragraphs by their text instead of a fixed index, so the\n# script keeps working even if paragraph numbering shifts slightly.\n$companyInfoIdx = 0\n$shortpositiesIdx = 0\n$shortPctIdx = 0\n$voorstelIdx = 0\n$overigeIdx = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($t -eq \"[COMPANY_INFO]`r\") { $companyInfoIdx = $i }\n    elseif ($t -eq \"Shortposities`r\") { $shortpositiesIdx = $i }\n    elseif ($t -eq \"[Short % of Shares Outstanding:]`r\") { $shortPctIdx = $i }\n    elseif ($t -eq \"Voorstel Koop & Verkoop`r\") { $voorstelIdx = $i }\n    elseif ($t -eq \"Overige`r\") { $overigeIdx = $i }\n}\n\n# 1) Add a new \"Logo:\" paragraph (style \"Plattetekst\") right before the\n#    \"[COMPANY_INFO]\" placeholder paragraph in the \"Bedrijfsinformatie\" section.\nif ($companyInfoIdx -gt 0) {\n    $companyInfoRange = $d.Paragraphs.Item($companyInfoIdx).Range\n    $companyInfoRange.InsertParagraphBefore()\n    $logoPara = $d.Paragraphs.Item($companyInfoIdx)\n    $logoPara.Range.Text = \"Logo:\"\n    $logoPara.Style = \"Plattetekst\"\n}\n\n# Re-resolve the remaining anchors: inserting the \"Logo:\" paragraph shifted\n# every paragraph index from $companyInfoIdx onward by +1.\nif ($shortpositiesIdx -ge $companyInfoIdx) { $shortpositiesIdx++ }\nif ($shortPctIdx -ge $companyInfoIdx) { $shortPctIdx++ }\nif ($voorstelIdx -ge $companyInfoIdx) { $voorstelIdx++ }\nif ($overigeIdx -ge $companyInfoIdx) { $overigeIdx++ }\n\n# 2) The <w:lastRenderedPageBreak/> marker moves from the first run of the\n#    \"[Short % of Shares Outstanding:]\" paragraph up onto the \"Shortposities\"\n#    heading run (it now renders at the top of that section instead).\nif ($shortpositiesIdx -gt 0) {\n    $r = $d.Paragraphs.Item($shortpositiesIdx).Range\n    $xml = Wrap-ParagraphPackage(\n        '<w:p><w:pPr><w:pStyle w:val=\"Kop4\"/><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr></w:pPr>' + `\n        '<w:bookmarkStart w:id=\"7\" w:name=\"shortposities\"/><w:bookmarkEnd w:id=\"6\"/>' + `\n        '<w:r><w:rPr><w:lang w:val=\"nl-NL\"/></w:rPr><w:lastRenderedPageBreak/><w:t>Shortposities</w:t></w:r></w:p>'\n    )\n    [void]$r.InsertXML($xml)\n}\n\nif ($shortPctIdx -gt 0) {\n    $r = $d.Paragraphs.Item($shortPctIdx).Range\n    $xml = Wrap-ParagraphPackage(\n        '<w:p><w:pPr><w:pStyle w:val=\"Plattetekst\"/></w:pPr>' + `\n        '<w:r><w:t>[</w:t></w:r><w:r><w:t>Short % of Shares Outstanding:</w:t></w:r><w:r><w:t>]</w:t></w:r></w:p>'\n    )\n    [void]$r.InsertXML($xml)\n}\n\n# 3) Add a <w:lastRenderedPageBreak/> to the \"Voorstel Koop & Verkoop\" heading.\nif ($voorstelIdx -gt 0) {\n    $r = $d.Paragraphs.Item($voorstelIdx).Range\n    $xml = Wrap-ParagraphPackage(\n        '<w:p><w:pPr><w:pStyle w:val=\"Kop4\"/></w:pPr>' + `\n        '<w:bookmarkStart w:id=\"14\" w:name=\"voorstel-koop-verkoop\"/><w:bookmarkEnd w:id=\"13\"/>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:lastRenderedPageBreak/><w:t>Voorstel</w:t></w:r><w:proofErr w:type=\"spellEnd\"/>' + `\n        '<w:r><w:t xml:space=\"preserve\"> Koop &amp; </w:t></w:r>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Verkoop</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n    )\n    [void]$r.InsertXML($xml)\n}\n\n# 4) Remove the <w:lastRenderedPageBreak/> from the \"Overige\" heading (the\n#    page break is now rendered earlier, at \"Voorstel Koop & Verkoop\").\nif ($overigeIdx -gt 0) {\n    $r = $d.Paragraphs.Item($overigeIdx).Range\n    $xml = Wrap-ParagraphPackage(\n        '<w:p><w:pPr><w:pStyle w:val=\"Kop4\"/></w:pPr>' + `\n        '<w:bookmarkStart w:id=\"15\" w:name=\"overige\"/><w:bookmarkEnd w:id=\"14\"/>' + `\n        '<w:proofErr w:type=\"spellStart\"/><w:r><w:t>Overige</w:t></w:r><w:proofErr w:type=\"spellEnd\"/></w:p>'\n    )\n    [void]$r.InsertXML($xml)\n}\n"}
